$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_cases_Smoke_learn.epam.com")

$ws.Range("A1:K1").Value = "Back on the main page"
$ws.Range("A2:K2").Value = "Test cases for Smoke Testing scope for learn.epam.com (on April 2023)"
